# AL_J5_jugadores.xlsx - per-matchday (jornada) reload of player stats.
# The incoming feed no longer carries a "Capitan" (captain) flag for this
# matchday, so column AP is cleared for every player row. Also normalizes
# the accented "Sí" spelling (previously unaccented "Si") in the
# "Suplente" column for the bench players (rows 13-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Capitan" column (AP) for every player row (2-21)
$ws.Range("AP2:AP21").ClearContents()

# Fix accent: "Si" -> "Sí" in the "Suplente" column (F) for rows 13-21
$ws.Range("F13").Value = "Sí"
$ws.Range("F14").Value = "Sí"
$ws.Range("F15").Value = "Sí"
$ws.Range("F16").Value = "Sí"
$ws.Range("F17").Value = "Sí"
$ws.Range("F18").Value = "Sí"
$ws.Range("F19").Value = "Sí"
$ws.Range("F20").Value = "Sí"
$ws.Range("F21").Value = "Sí"
